# String to int[] 함수 추가
# The "weatherType" column (D) used bracketed JSON-style arrays like
# "[0,1,2,3,4,5]" as the encoded default value. Switch it to a
# semicolon-delimited form "0;1;2;3;4;5" that the new String->int[]
# parser expects.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2..5) {
    $ws.Cells.Item($r, 4).Value = "0;1;2;3;4;5"
}

# Leave the active selection on D2, matching the saved workbook state.
$ws.Range("D2").Select()
